$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty second row (row 2), shifting all data rows up by one.
$ws.Rows.Item(2).Delete()

# Update the selection to match the new active cell/selection (entire row 2 selected).
$ws.Range("A2:XFD2").Select()
